$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Generated Features")

$ws.Range("A4").Value = "groc14_per_1k_capita"
$ws.Range("A5").Value = "superc14_per_1k_capita"
$ws.Range("A6").Value = "convs14_per_1k_capita"
$ws.Range("A7").Value = "specs14_per_1k_capita"

$ws.Range("B4").Value = "Count of grocery stores in county 2014 per 1,000 capita."
$ws.Range("B5").Value = "Count of supercenter stores in county 2014 per 1,000 capita."
$ws.Range("B6").Value = "Count of convenience stores in county 2014 per 1,000 capita."
$ws.Range("B7").Value = "Count of specialty food stores in county 2014 per 1,000 capita."
